$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at J:K. This pushes the old J/K/L columns
# (CHECK_GROUND_TRUTH / GROUND_TRUTH_PATH / DATABASE_PATH) to L/M/N.
$ws.Range("J:K").Insert()

# The two freshly inserted columns should share column I's width (14 chars),
# matching the merged <col min="9" max="11" .../> range in the target file.
$iWidth = $ws.Range("I:I").ColumnWidth
$ws.Range("J:K").ColumnWidth = $iWidth

# New header text for the inserted columns
$ws.Range("J1").Value = "REGION_TABLE_PATH"
$ws.Range("K1").Value = "REFER_TEXT_TABLE_PATH"

# Update the active selection to match the target workbook
[void]$ws.Range("O18").Select()
